# CARAGA_REPAIR.xlsx formatting fix:
#   - drop the stale "PREVIOUS ACCOMPLISHMENT" .. "Unnamed: 46" columns (AM:AU),
#     which shifts the trailing "Status as of July 4, 2025" column back to AM
#   - backfill the new site/classroom-status breakdown columns (AA:AL) for each
#     data row, and make sure every touched cell (header + body) carries a
#     thin border so the sheet reads as a proper bordered table
#   - the status column header (AM1) keeps the workbook's existing bold header
#     look, the dropdown validation list follows the column to AM2:AM5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 5

# ---------------------------------------------------------------------------
# 1) Remove the obsolete AM:AU columns. This deletes the old
#    "PREVIOUS ACCOMPLISHMENT", "DIFFERENCE", "Projected Date of Completion",
#    "Month Completed", "No. of Sites (original)", "No. of Classrooms
#    (original)", "Project Allocation (original)", "BBM TERM COMPLETED" and
#    "Unnamed: 46" columns in one shot, and shifts the last column ("Status
#    as of July 4, 2025", formerly AV) left into AM - carrying the
#    worksheet's dimension, data validation (AV2:AV5 -> AM2:AM5) along with
#    it automatically.
$ws.Range("AM1:AU5").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2) Backfill the site/classroom status-breakdown columns (AA:AL) for every
#    data row. These projects are all COMPLETED, so every site/classroom
#    counts under the "Completed" buckets (AE / AJ) and every other bucket is
#    zero.
for ($r = 2; $r -le $lastRow; $r++) {
    $totalSites = $ws.Cells.Item($r, 8).Value2   # column H - TOTAL NO. OF SITES
    $totalRooms = $ws.Cells.Item($r, 9).Value2   # column I - TOTAL PHYSICAL TARGET

    $ws.Cells.Item($r, 27).Value = 0            # AA - No. of Sites Reverted
    $ws.Cells.Item($r, 28).Value = 0            # AB - No. of Sites Not yet started
    $ws.Cells.Item($r, 29).Value = 0            # AC - No. of Sites Under Procurement
    $ws.Cells.Item($r, 30).Value = 0            # AD - No. of Sites On Going
    $ws.Cells.Item($r, 31).Value = $totalSites  # AE - No. of Sites Completed
    $ws.Cells.Item($r, 32).Value = 0            # AF - No. of CL Reverted
    $ws.Cells.Item($r, 33).Value = 0            # AG - No. of CL Not yet started
    $ws.Cells.Item($r, 34).Value = 0            # AH - No. of CL Under Procurement
    $ws.Cells.Item($r, 35).Value = 0            # AI - No. of CL On Going
    $ws.Cells.Item($r, 36).Value = $totalRooms  # AJ - No. of CL Completed
    $ws.Cells.Item($r, 37).Value = 0            # AK - No. of Sites Terminated
    $ws.Cells.Item($r, 38).Value = 0            # AL - No. of CL Terminated
}

# ---------------------------------------------------------------------------
# 3) Give every used cell - header and body alike - a thin border, so the
#    whole A1:AM5 table is boxed in. Doing this first (while AM1 is still
#    plain) lets the border reuse the workbook's existing thin-border
#    definition instead of minting a duplicate one.
$usedRange = $ws.Range("A1:AM" + $lastRow)
$usedRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 4) AM1 ("Status as of July 4, 2025") should look like the rest of the bold
#    header row, but without the header's centered/top alignment. Copy A1's
#    format (bold font + border) onto it, then reset alignment back to the
#    sheet default so only the bold weight carries over.
$ws.Range("A1").Copy()
$ws.Range("AM1").PasteSpecial(-4122)
$ws.Range("AM1").HorizontalAlignment = 1
$ws.Range("AM1").VerticalAlignment = -4107
$excel.CutCopyMode = $false

Write-Output ("Final used range: " + $ws.UsedRange.Address())
